$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item catalogue for rows 2-6: row 2 replaces the old sample item
# (previously "SQ009" / "SQ009-Pulsera") and rows 3-6 are brand-new rows
# appended below it.
$names  = @("SP001-Polera", "SS016-Plancha", "SZ006-ZAPATERA", "SA042-RIZADOR", "SG015-ZAPATO")
$codes  = @("SP001",        "SS016",         "SZ006",          "SA042",         "SG015")
$prices = @(80,              60,              20,               80,              30)

# Write column by column (Nombre, then Codigo Interno, then the rest) so
# the data lands the same way it would from a column-oriented paste.
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $names[$i]
}
for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = $i + 2
    $ws.Range("B$r").Value = $codes[$i]
}

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Range("E$r").Value = "NIU"
    $ws.Range("F$r").Value = "PEN"
    $ws.Range("G$r").Value = $prices[$i]
    $ws.Range("H$r").Value = 10
    $ws.Range("I$r").Value = "SI"
    $ws.Range("K$r").Value = 10
    $ws.Range("L$r").Value = 1
    $ws.Range("M$r").Value = 0
    $ws.Range("P$r").Value = $names[$i]
    $ws.Range("T$r").Value = $codes[$i]
}
